# final_project/db_design.xlsx — "add final project, project"
#
# Adds three new sheets (Grades comments, Suspended students, Enrollments),
# renumbers the Courses index column from 0-based to 1-based, and points the
# sole Grades row at the real student ID. Also restores the various sheet
# selections to their final (post-edit) state.

$wb = $excel.ActiveWorkbook

$grades      = $wb.Worksheets.Item("Grades")
$instructors = $wb.Worksheets.Item("Instructors")
$courses     = $wb.Worksheets.Item("Courses")
$students    = $wb.Worksheets.Item("Students")
$fields      = $wb.Worksheets.Item("Fields of study")

# ---------------------------------------------------------------------
# 1. New sheet: Enrollments (inserted right after Instructors).
#    Typed first so the shared string "Student_Index" is introduced
#    before "Comment"/"Grade_ID" (matches the shared-strings table order
#    produced by the original authoring session).
# ---------------------------------------------------------------------
$enrollments = $wb.Worksheets.Add($null, $instructors)
$enrollments.Name = "Enrollments"
$enrollments.Range("A1").Value = "Index"
$enrollments.Range("B1").Value = "Course_ID"
$enrollments.Range("C1").Value = "Student_Index"
$enrollments.Range("A2").Value = 1
$enrollments.Range("B2").Value = 1
$enrollments.Range("C2").Value = 100101

# ---------------------------------------------------------------------
# 2. New sheet: Grades comments (inserted right after Grades).
# ---------------------------------------------------------------------
$gradesComments = $wb.Worksheets.Add($null, $grades)
$gradesComments.Name = "Grades comments"
$gradesComments.Range("B1").Value = "Comment"
$gradesComments.Range("A1").Value = "Grade_ID"
$gradesComments.Range("A2").Value = 1
$gradesComments.Range("B2").Value = "Some comment."

# ---------------------------------------------------------------------
# 3. New sheet: Suspended students (inserted right after Grades comments).
# ---------------------------------------------------------------------
$suspended = $wb.Worksheets.Add($null, $gradesComments)
$suspended.Name = "Suspended students"
$suspended.Range("A1").Value = "Student_Index"
$suspended.Range("B1").Value = "Reason"
$suspended.Range("A2").Value = 100102
$suspended.Range("B2").Value = "Some reason…"

# ---------------------------------------------------------------------
# 4. Data edits on existing sheets.
# ---------------------------------------------------------------------

# Courses: the Index column moves from 0-based to 1-based numbering.
$courses.Range("A2").Value = 1
$courses.Range("A3").Value = 2
$courses.Range("A4").Value = 3
$courses.Range("A5").Value = 4
$courses.Range("A6").Value = 5

# Grades: the lone row now references the real Student_ID.
$grades.Range("C2").Value = 100101

# ---------------------------------------------------------------------
# 5. Selections (sheetView state) for each sheet, left-to-right, ending
#    on Students so it is the active tab on reopen.
# ---------------------------------------------------------------------

$courses.Range("A7").Select()

$fields.Range("C31").Select()

$grades.Range("L4").Select()

$gradesComments.Range("B3").Select()

$suspended.Columns("B").Select()

$instructors.Range("C31:C32").Select()

$enrollments.Range("C3").Select()

$students.Activate()
$students.Range("S5").Select()
